$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run NATMI Slit2-Robo2 results now include the "ECs" cluster, giving a full 3x3
# Sending-cluster x Target-cluster matrix (ECs / FAPs / sCs) instead of the previous 2x2.
$nRows = 9
$nCols = 20
$data = New-Object 'object[,]' $nRows,$nCols

$data[0,0] = "ECs"
$data[0,1] = "Slit2"
$data[0,2] = "Robo2"
$data[0,3] = "ECs"
$data[0,4] = 2
$data[0,5] = 0.6666666666666666
$data[0,6] = 0.143896
$data[0,7] = 0.431688
$data[0,8] = 0.02807111181859822
$data[0,9] = 0.02807111181859822
$data[0,10] = 2
$data[0,11] = 0.6666666666666666
$data[0,12] = 1.054425
$data[0,13] = 3.163275
$data[0,14] = 0.5590543736060202
$data[0,15] = 0.5590543736060202
$data[0,16] = 0.1517275398
$data[0,17] = 1.3655478582
$data[0,18] = 0.01569327783417098
$data[0,19] = 0.01569327783417098

$data[1,0] = "ECs"
$data[1,1] = "Slit2"
$data[1,2] = "Robo2"
$data[1,3] = "FAPs"
$data[1,4] = 2
$data[1,5] = 0.6666666666666666
$data[1,6] = 0.143896
$data[1,7] = 0.431688
$data[1,8] = 0.02807111181859822
$data[1,9] = 0.02807111181859822
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 0.7756663333333332
$data[1,13] = 2.326999
$data[1,14] = 0.4112569941996302
$data[1,15] = 0.4112569941996303
$data[1,16] = 0.1116152827013333
$data[1,17] = 1.004537544312
$data[1,18] = 0.01154444107035842
$data[1,19] = 0.01154444107035842

$data[2,0] = "ECs"
$data[2,1] = "Slit2"
$data[2,2] = "Robo2"
$data[2,3] = "sCs"
$data[2,4] = 2
$data[2,5] = 0.6666666666666666
$data[2,6] = 0.143896
$data[2,7] = 0.431688
$data[2,8] = 0.02807111181859822
$data[2,9] = 0.02807111181859822
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.05599533333333334
$data[2,13] = 0.167986
$data[2,14] = 0.0296886321943495
$data[2,15] = 0.02968863219434951
$data[2,16] = 0.008057504485333334
$data[2,17] = 0.072517540368
$data[2,18] = 0.0008333929140688201
$data[2,19] = 0.0008333929140688202

$data[3,0] = "FAPs"
$data[3,1] = "Slit2"
$data[3,2] = "Robo2"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 4.277274333333334
$data[3,7] = 12.831823
$data[3,8] = 0.8344071140950421
$data[3,9] = 0.8344071140950421
$data[3,10] = 2
$data[3,11] = 0.6666666666666666
$data[3,12] = 1.054425
$data[3,13] = 3.163275
$data[3,14] = 0.5590543736060202
$data[3,15] = 0.5590543736060202
$data[3,16] = 4.510064988925
$data[3,17] = 40.590584900325
$data[3,18] = 0.4664789465028107
$data[3,19] = 0.4664789465028107

$data[4,0] = "FAPs"
$data[4,1] = "Slit2"
$data[4,2] = "Robo2"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 4.277274333333334
$data[4,7] = 12.831823
$data[4,8] = 0.8344071140950421
$data[4,9] = 0.8344071140950421
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 0.7756663333333332
$data[4,13] = 2.326999
$data[4,14] = 0.4112569941996302
$data[4,15] = 0.4112569941996303
$data[4,16] = 3.317737698797444
$data[4,17] = 29.859639289177
$data[4,18] = 0.343155761681515
$data[4,19] = 0.343155761681515

$data[5,0] = "FAPs"
$data[5,1] = "Slit2"
$data[5,2] = "Robo2"
$data[5,3] = "sCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 4.277274333333334
$data[5,7] = 12.831823
$data[5,8] = 0.8344071140950421
$data[5,9] = 0.8344071140950421
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.05599533333333334
$data[5,13] = 0.167986
$data[5,14] = 0.0296886321943495
$data[5,15] = 0.02968863219434951
$data[5,16] = 0.2395074020531112
$data[5,17] = 2.155566618478
$data[5,18] = 0.02477240591071633
$data[5,19] = 0.02477240591071633

$data[6,0] = "sCs"
$data[6,1] = "Slit2"
$data[6,2] = "Robo2"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 0.7049536666666666
$data[6,7] = 2.114861
$data[6,8] = 0.1375217740863597
$data[6,9] = 0.1375217740863597
$data[6,10] = 2
$data[6,11] = 0.6666666666666666
$data[6,12] = 1.054425
$data[6,13] = 3.163275
$data[6,14] = 0.5590543736060202
$data[6,15] = 0.5590543736060202
$data[6,16] = 0.7433207699749999
$data[6,17] = 6.689886929774999
$data[6,18] = 0.07688214926903844
$data[6,19] = 0.07688214926903846

$data[7,0] = "sCs"
$data[7,1] = "Slit2"
$data[7,2] = "Robo2"
$data[7,3] = "FAPs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 0.7049536666666666
$data[7,7] = 2.114861
$data[7,8] = 0.1375217740863597
$data[7,9] = 0.1375217740863597
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 0.7756663333333332
$data[7,13] = 2.326999
$data[7,14] = 0.4112569941996302
$data[7,15] = 0.4112569941996303
$data[7,16] = 0.5468088257932221
$data[7,17] = 4.921279432138999
$data[7,18] = 0.0565567914477569
$data[7,19] = 0.05655679144775692

$data[8,0] = "sCs"
$data[8,1] = "Slit2"
$data[8,2] = "Robo2"
$data[8,3] = "sCs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 0.7049536666666666
$data[8,7] = 2.114861
$data[8,8] = 0.1375217740863597
$data[8,9] = 0.1375217740863597
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.05599533333333334
$data[8,13] = 0.167986
$data[8,14] = 0.0296886321943495
$data[8,15] = 0.02968863219434951
$data[8,16] = 0.03947411554955556
$data[8,17] = 0.355267039946
$data[8,18] = 0.004082833369564358
$data[8,19] = 0.004082833369564359

$startRow = 2
$endRow = $startRow + $nRows - 1
$rng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, $nCols))
$rng.Value = $data

Write-Output "done"